# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
#
# Updates the Shock (sTop/xMax/xMin) hardpoint coordinates on both the
# front ("S2LAR_Sedan_HambaLG_f") and rear ("S2LAR_Sedan_HambaLG_r")
# sheets.

$wb = $excel.ActiveWorkbook

$wsFront = $wb.Worksheets.Item("S2LAR_Sedan_HambaLG_f")
$wsRear  = $wb.Worksheets.Item("S2LAR_Sedan_HambaLG_r")

# ---------------------------------------------------------------------
# Front sheet (S2LAR_Sedan_HambaLG_f)
# ---------------------------------------------------------------------

# Row 23 - Shock / sTop
$wsFront.Range("F23").Value = 0.15379999999999999
$wsFront.Range("G23").Value = 0.65
$wsFront.Range("H23").Value = 0.24

# Row 24 - Shock / sBottom
$wsFront.Range("G24").Value = 0.91
$wsFront.Range("H24").Value = 0.23

# Row 26 - Endstop / xMax (drop the extra decimal place of precision,
# matching the format already used by row 23/24)
$wsFront.Range("G26").NumberFormat = "0.00"
$wsFront.Range("H26").NumberFormat = "0.00"
$wsFront.Range("G26").Value = 0.62
$wsFront.Range("H26").Value = 0.65

# Row 27 - Endstop / xMin
$wsFront.Range("G27").NumberFormat = "0.00"
$wsFront.Range("H27").NumberFormat = "0.00"
$wsFront.Range("G27").Value = 0.85
$wsFront.Range("H27").Value = 0.19

# ---------------------------------------------------------------------
# Rear sheet (S2LAR_Sedan_HambaLG_r)
# ---------------------------------------------------------------------

# Row 23 - Shock / sTop
$wsRear.Range("F23").Value = 0.13
$wsRear.Range("G23").Value = 0.65
$wsRear.Range("H23").Value = 0.24

# Rows 26/27 column F hold the same hardpoint numbers as before, just
# re-settled by a hair of floating-point noise after the recalculation
# that produced this revision.
$wsRear.Range("F26").Value = 2.6557142857142869E-3
$wsRear.Range("F27").Value = -5.5166428571428582E-2

# Row 24 - Shock / sBottom
$wsRear.Range("F24").Value = 0.13
$wsRear.Range("G24").Value = 0.91
$wsRear.Range("H24").Value = 0.23

# Row 26 - Endstop / xMax
$wsRear.Range("G26").NumberFormat = "0.00"
$wsRear.Range("H26").NumberFormat = "0.00"
$wsRear.Range("G26").Value = 0.62
$wsRear.Range("H26").Value = 0.65

# Row 27 - Endstop / xMin
$wsRear.Range("G27").NumberFormat = "0.00"
$wsRear.Range("H27").NumberFormat = "0.00"
$wsRear.Range("G27").Value = 0.85
$wsRear.Range("H27").Value = 0.19

$excel.Calculate()
